# Резюме.docx edit script
# Applies the changes described by the commit:
#  1. Move the automatic "_GoBack" bookmark from the paragraph right after
#     the "Резюме" heading to the (empty) paragraph right after the big
#     summary paragraph - this is exactly what Word does automatically
#     when the last edit location changes.
#  2. Insert a missing space between "използваните" and "(скрити модели".
#  3. Unlink the "REF _Ref341289950" cross-reference field so that the
#     title ("Допълнителни имплементации и тестове") is left behind as
#     plain text instead of a field.
#  4. Split "ченическа конференция" into "ченическа " / "секция" in the
#     page header (i.e. rename "конференция" to "секция").

$d = $word.ActiveDocument

# --- 1. Insert the missing space before "(скрити модели" -------------------
$d.Content.Find.Execute(
    "използваните(скрити",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "използваните (скрити", 2) | Out-Null

# --- 2. Turn the REF field into plain text ---------------------------------
if ($d.Fields.Count -ge 1) {
    $d.Fields(1).Unlink()
}

# --- 3. Move the "_GoBack" bookmark -----------------------------------------
# Find the paragraph that follows the one ending in
# "... зависимости. " (the empty paragraph right after the big summary
# paragraph) and (re)plant the bookmark there. Adding a bookmark with a
# name that already exists moves it, exactly like Word's own "_GoBack"
# bookkeeping.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "зависимости") {
        $target = $d.Paragraphs($i + 1).Range
        $d.Bookmarks.Add("_GoBack", $target) | Out-Null
        break
    }
}

# --- 4. Header: "ченическа конференция" -> "ченическа секция" --------------
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Find.Execute(
    "конференция",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "секция", 2) | Out-Null
